$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1200.4
$ws.Range("I18").Value = 1375.5
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 1375.5
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = -1091.5
$ws.Range("N18").Value = -1068

$ws.Range("H21").Value = 29990
$ws.Range("I21").Value = 29950
$ws.Range("K21").Value = 29950
$ws.Range("M21").Value = -29482

$ws.Range("H23").Value = 29990
$ws.Range("I23").Value = 29950
$ws.Range("K23").Value = 29950
$ws.Range("M23").Value = -29716

$ws.Range("H113").Value = 3333.3333
$ws.Range("I113").Value = 2500
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 2500
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = 754
$ws.Range("N113").Value = -10008

$ws.Range("H137").Value = 1215.0819
$ws.Range("I137").Value = 1009.1556
$ws.Range("J137").Value = 1794.25
$ws.Range("K137").Value = 3027.4668
$ws.Range("L137").Value = 5382.75
$ws.Range("M137").Value = -477.4668000000001
$ws.Range("N137").Value = -10482.75

$ws.Range("H138").Value = 2764.243
$ws.Range("I138").Value = 1797.875
$ws.Range("K138").Value = 5393.625
$ws.Range("M138").Value = -253.625

$ws.Range("H140").Value = 96248.336
$ws.Range("J140").Value = 107998
$ws.Range("L140").Value = 107998
$ws.Range("N140").Value = -118358

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12984.777
$ws.Range("I32").Value = 13821.267
$ws.Range("J32").Value = 8802.333000000001
$ws.Range("K32").Value = 13821.267
$ws.Range("L32").Value = 8802.333000000001
$ws.Range("M32").Value = -13534.267
$ws.Range("N32").Value = -9376.333000000001

$ws.Range("H45").Value = 1730.6364
$ws.Range("I45").Value = 1420.6666
$ws.Range("J45").Value = 2102.6
$ws.Range("K45").Value = 1420.6666
$ws.Range("L45").Value = 2102.6
$ws.Range("M45").Value = -1043.6666
$ws.Range("N45").Value = -2856.6

$ws.Range("H61").Value = 967.0536
$ws.Range("I61").Value = 946.06525
$ws.Range("J61").Value = 1063.6
$ws.Range("K61").Value = 946.06525
$ws.Range("L61").Value = 1063.6
$ws.Range("M61").Value = -734.06525
$ws.Range("N61").Value = -1487.6

$ws.Range("H110").Value = 1878.1818
$ws.Range("I110").Value = 1866.1
$ws.Range("K110").Value = 1866.1
$ws.Range("M110").Value = 178.9000000000001

$ws.Range("H132").Value = 2934.138
$ws.Range("I132").Value = 2442.25
$ws.Range("J132").Value = 3539.5386
$ws.Range("K132").Value = 7326.75
$ws.Range("L132").Value = 10618.6158
$ws.Range("M132").Value = -4796.75
$ws.Range("N132").Value = -15678.6158

$ws.Range("H136").Value = 967.0536
$ws.Range("I136").Value = 946.06525
$ws.Range("J136").Value = 1063.6
$ws.Range("K136").Value = 2838.19575
$ws.Range("L136").Value = 3190.8
$ws.Range("M136").Value = -288.1957499999999
$ws.Range("N136").Value = -8290.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 15235
$ws.Range("I5").Value = 9369
$ws.Range("J5").Value = 29900
$ws.Range("K5").Value = 9369
$ws.Range("L5").Value = 29900
$ws.Range("M5").Value = -9256
$ws.Range("N5").Value = -30126

$ws.Range("H22").Value = 25500
$ws.Range("I22").Value = 50000
$ws.Range("K22").Value = 50000
$ws.Range("M22").Value = -49827

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5559.263
$ws.Range("I22").Value = 7919.5386
$ws.Range("K22").Value = 7919.5386
$ws.Range("M22").Value = -7569.5386

$ws.Range("H31").Value = 2092.7646
$ws.Range("I31").Value = 1378.7307
$ws.Range("J31").Value = 4413.375
$ws.Range("K31").Value = 1378.7307
$ws.Range("L31").Value = 4413.375
$ws.Range("M31").Value = -1083.7307
$ws.Range("N31").Value = -5003.375

$ws.Range("H34").Value = 2092.7646
$ws.Range("I34").Value = 1378.7307
$ws.Range("J34").Value = 4413.375
$ws.Range("K34").Value = 1378.7307
$ws.Range("L34").Value = 4413.375
$ws.Range("M34").Value = -1176.7307
$ws.Range("N34").Value = -4817.375

$ws.Range("H134").Value = 1567.3636
$ws.Range("I134").Value = 1280.6923
$ws.Range("J134").Value = 2632.1428
$ws.Range("K134").Value = 3842.0769
$ws.Range("L134").Value = 7896.428400000001
$ws.Range("M134").Value = -1307.0769
$ws.Range("N134").Value = -12966.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 173.13333
$ws.Range("I7").Value = 111.7
$ws.Range("J7").Value = 296
$ws.Range("K7").Value = 335.1
$ws.Range("L7").Value = 888
$ws.Range("M7").Value = -223.1
$ws.Range("N7").Value = -1112

$ws.Range("H33").Value = 975.1429000000001
$ws.Range("J33").Value = 2571
$ws.Range("L33").Value = 15426
$ws.Range("N33").Value = -15992

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 4898.9546
$ws.Range("I107").Value = 7255.7856
$ws.Range("J107").Value = 774.5
$ws.Range("K107").Value = 7255.7856
$ws.Range("L107").Value = 774.5
$ws.Range("M107").Value = -5335.7856
$ws.Range("N107").Value = -4614.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 10108.571
$ws.Range("I100").Value = 18260
$ws.Range("K100").Value = 18260
$ws.Range("M100").Value = -17719

$ws.Range("H122").Value = 20459710
$ws.Range("I122").Value = 19236062
$ws.Range("J122").Value = 22227200
$ws.Range("K122").Value = 57708186
$ws.Range("L122").Value = 66681600
$ws.Range("M122").Value = -57705736
$ws.Range("N122").Value = -66686500

$ws.Range("H136").Value = 1549
$ws.Range("I136").Value = 1408.1549
$ws.Range("J136").Value = 2049
$ws.Range("K136").Value = 4224.4647
$ws.Range("L136").Value = 6147
$ws.Range("M136").Value = -1674.4647
$ws.Range("N136").Value = -11247

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 24042378
$ws.Range("I122").Value = 50003600
$ws.Range("J122").Value = 7816614
$ws.Range("K122").Value = 150010800
$ws.Range("L122").Value = 23449842
$ws.Range("M122").Value = -150008350
$ws.Range("N122").Value = -23454742

$ws.Range("H141").Value = 22857.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 22857.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 22857.5
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -33217.5
